# LOB1208.xlsx content fix-up
# The sheet had several rows whose B/C (data) columns were off-by-one
# relative to their A column (label): values had "leaked" from the row
# below/above. This script:
#   1. Fixes row 10 (Objetivos:) which held Robson's name instead of the
#      Portuguese "Objetivos" text.
#   2. Inserts a new row 13 to hold "Docentes responsáveis:" data
#      (Robson's name), shifting everything below down by one row.
#   3. Fixes the (now shifted) rows whose data text no longer matched
#      their label, restoring/adding the correct Portuguese content for
#      "Programa resumido:", "Programa:", "Método:", "Critério:",
#      "Norma de recuperação:" and "Bibliografia:".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Row 10 "Objetivos:" - replace wrong content with the real PT text
# ---------------------------------------------------------------------
$objetivosPT = @"
Formar profissionais em nível superior com capacidade de conhecer a sequência dos procedimentos de análise química de interesse ambiental desde a obtenção das amostras in situ até o preparo preliminar do material a ser analisado. Executar procedimentos de análises qualitativas bem como interpretar, avaliar e criticar os resultados obtidos. Objetivos Específicos: Formar profissionais em nível superior com capacidade de conhecer as etapas da sequência analítica. Compreender e aplicar os procedimentos mais comuns de amostragem, coleta e preparação de amostras bem como os erros a não cometer no preparo das amostras de interesse ambiental. Empregar tratamentos preliminares no preparo das amostras: limpeza, secagem, moagem e peneiramento. Compreender as bases teóricas da química analítica qualitativa de interesse ambiental.
"@
$objetivosPT = $objetivosPT.Trim()

$ws.Range("B10").Value = $objetivosPT
$ws.Range("C10").Value = $objetivosPT

# ---------------------------------------------------------------------
# 2) Insert a new row at 13 for "Docentes responsáveis:" data, pushing
#    the old rows 13-23 down to 14-24. Copy B/C formatting from the row
#    that lands right below (the new row 14) so styles match the rest
#    of the table, then clear the (empty) A13 cell.
# ---------------------------------------------------------------------
$ws.Rows(13).Insert()

$ws.Range("A14:C14").Copy()
$ws.Range("A13:C13").PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = $false
$ws.Range("A13").Clear()

$robson = "7455355 - Robson da Silva Rocha"
$ws.Range("B13").Value = $robson
$ws.Range("C13").Value = $robson

# ---------------------------------------------------------------------
# 3) Fix the shifted rows' B/C content
# ---------------------------------------------------------------------

# Row 14 "Programa resumido:" -> Portuguese short syllabus
$programaResumido = @"
Introdução à análise qualitativa, indicando suas aplicabilidades e limitações. Uso das técnicas qualitativas para análise dos principais íons de importância ambiental. Análise de sólidos, partículas, sedimentos. Estudos de amostras de importância ambiental.
"@
$programaResumido = $programaResumido.Trim()
$ws.Range("B14").Value = $programaResumido
$ws.Range("C14").Value = $programaResumido

# Row 16 "Programa:" -> Portuguese full syllabus
$programa = @"
- Revisão das regras de segurança laboratorial - Introdução à análise qualitativa: Definições, objetivos e limitações. - Análise de sólidos, partículas, sedimentos.- Identificação dos cátions do grupo I (K+, Na+ e NH4+); grupo II (Mg2+, Ca2+ e Ba2+); grupo III (Al3+, Fe3+, Mn2+).- Estudo dos ânions e suas aplicações em análises ambientais (Cl e suas espécies, SO42-, CO32-, S2-, NO3-).- Análise gravimétrica: fundamentos e formação de precipitados.- Análises dos principais cátions e ânions em amostras conhecidas e desconhecidas para os alunos- Análise de metais em solo, água ou outras amostras ambientais importantes
"@
$programa = $programa.Trim()
$ws.Range("B16").Value = $programa
$ws.Range("C16").Value = $programa

# Row 19 "Método:" -> method text
$metodo = @"
O método de avaliação será composto por avaliações teóricas, práticas e relatórios de atividades de práticas laboratoriais.
"@
$metodo = $metodo.Trim()
$ws.Range("B19").Value = $metodo
$ws.Range("C19").Value = $metodo

# Row 20 "Critério:" -> criteria text
$criterio = @"
Para o cálculo da nota final (NF) será feita a média aritmética das avaliações aplicadas. Estará aprovado por notas o aluno que obtiver nota final igual ou superior a 5,0 pontos.
"@
$criterio = $criterio.Trim()
$ws.Range("B20").Value = $criterio
$ws.Range("C20").Value = $criterio

# Row 21 "Norma de recuperação:" -> recovery norm text
$norma = "Avaliação de recuperação (R) envolvendo todo o conteúdo da disciplina. Média Final = (NF+R) / 2 => 5,0 Aprovado"
$ws.Range("B21").Value = $norma
$ws.Range("C21").Value = $norma

# Row 22 "Bibliografia:" -> bibliography text (brand new content)
$bibliografia = @"
Baird, C., Michael C.C.  Environmental chemistry. Editora: New York: Freeman, 5a edição, 2012Baird, C.,Michael C.C.  Química ambiental. Editora: Porto Alegre: Bookman,  4a edição, 2011Harris, D. C. Análise Química Quantitativa. Editora: LTC, 8a edição, 2012Krug, F.J., Rocha F.R.P. Métodos de preparo de amostras para análise elementar. Editora EditSBQ, 1a edição, 2016Luna, A. Química analítica ambiental. Editora: EdUERJ, 1a edição, 2003Rocha, J.C., Rosa, A.H., Cardoso, A.A. Introdução à química ambiental. Editora: Porto Alegre: Bookman, 2a edição, 2009.Skoog, D. A, West, D. M., Holler, F. J., Crouch, S. R. Fundamentos de Química Analítica. Editora: Thomson, tradução da 8ª edição, 2006
"@
$bibliografia = $bibliografia.Trim()
$ws.Range("B22").Value = $bibliografia
$ws.Range("C22").Value = $bibliografia

# ---------------------------------------------------------------------
# 4) Column widths: column A keeps its own 30.71 width/style, column B
#    keeps 60.71 (these were previously lumped into one "1-2" col
#    entry sharing A's width/style; split them so B has its own entry).
# ---------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 30.7109375
$ws.Columns("B").ColumnWidth = 60.7109375
